$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19014048576355
$ws.Range("B1").Value = 1.847458839416504
$ws.Range("C1").Value = 6.666423320770264
$ws.Range("D1").Value = 2.284687519073486
$ws.Range("E1").Value = 1.193082451820374
